$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-22 Tuesday", "2025-04-23 Wednesday"),
    @("905×6=5430", "698×3=2094"),
    @("351×6=2106", "209×4=836"),
    @("966×7=6762", "114×4=456"),
    @("714×8=5712", "437×2=874"),
    @("326×9=2934", "369×6=2214"),
    @("513×9=4617", "543×8=4344"),
    @("156×9=1404", "204×6=1224"),
    @("593×2=1186", "178×2=356"),
    @("687×7=4809", "910×2=1820"),
    @("887×2=1774", "400×3=1200"),
    @("602×2=1204", "349×3=1047"),
    @("941×4=3764", "206×2=412"),
    @("798×3=2394", "454×6=2724"),
    @("816×7=5712", "923×9=8307"),
    @("563×5=2815", "353×2=706"),
    @("364×4=1456", "315×6=1890"),
    @("172×5=860", "450×4=1800"),
    @("679×3=2037", "418×3=1254"),
    @("399×7=2793", "254×5=1270"),
    @("983×4=3932", "981×7=6867"),
    @("739×9=6651", "109×8=872"),
    @("629×8=5032", "230×6=1380"),
    @("599×9=5391", "996×9=8964"),
    @("925×9=8325", "678×9=6102"),
    @("379×6=2274", "690×3=2070")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
